$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("diem_danh")

# Header row
$ws.Range("A1").Value = "buoi_hoc_id"
$ws.Range("B1").Value = "hoc_sinh_id"
$ws.Range("C1").Value = "ketqua"

# Data rows: buoi_hoc_id, hoc_sinh_id, ketqua
$data = @(
    @(1,1,0),
    @(1,2,1),
    @(1,3,-1),
    @(1,4,0),
    @(2,1,0),
    @(2,2,0),
    @(2,3,0),
    @(2,4,0),
    @(3,1,0),
    @(3,2,0),
    @(3,3,0),
    @(3,4,0),
    @(4,1,0),
    @(4,2,0),
    @(4,3,0),
    @(4,4,0)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Column widths (closest attainable values via the ColumnWidth property,
# which Excel snaps to the underlying pixel grid)
$ws.Columns.Item(1).ColumnWidth = 12.65
$ws.Columns.Item(2).ColumnWidth = 12.8

# Page setup (portrait)
$ws.PageSetup.Orientation = 1

# Make this sheet the active one, with the given selection
$ws.Activate()
$ws.Range("E18").Select()
